$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.059224843978882
$ws.Range("B1").Value = 1.254661560058594
$ws.Range("C1").Value = 1.664891481399536
$ws.Range("D1").Value = 3.978271484375
$ws.Range("E1").Value = 3.413884162902832
